# Update NATMI LR-pair metrics (Mfge8-Itgav) with recomputed TPM-derived values.
# Only the "ECs" sending/target cluster's ligand (G/H) and receptor (M/N)
# expression values changed; all downstream derived-specificity and edge
# weight columns (I,J,O,P,Q,R,S,T) are recomputed accordingly for every row
# that references the ECs cluster as sender and/or target.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 21.64449733333333
$ws.Range("H2").Value2 = 64.933492
$ws.Range("I2").Value2 = 0.1098676276771345
$ws.Range("J2").Value2 = 0.1098676276771345
$ws.Range("M2").Value2 = 3.759736666666667
$ws.Range("N2").Value2 = 11.27921
$ws.Range("O2").Value2 = 0.0683751702595819
$ws.Range("P2").Value2 = 0.06837517025958188
$ws.Range("Q2").Value2 = 81.37761025570222
$ws.Range("R2").Value2 = 732.39849230132
$ws.Range("S2").Value2 = 0.007512217748440423
$ws.Range("T2").Value2 = 0.007512217748440423
$ws.Range("G3").Value2 = 21.64449733333333
$ws.Range("H3").Value2 = 64.933492
$ws.Range("I3").Value2 = 0.1098676276771345
$ws.Range("J3").Value2 = 0.1098676276771345
$ws.Range("O3").Value2 = 0.6514180024294648
$ws.Range("P3").Value2 = 0.6514180024294647
$ws.Range("Q3").Value2 = 775.2937230576662
$ws.Range("R3").Value2 = 6977.643507518997
$ws.Range("S3").Value2 = 0.07156975055310313
$ws.Range("T3").Value2 = 0.07156975055310313
$ws.Range("G4").Value2 = 21.64449733333333
$ws.Range("H4").Value2 = 64.933492
$ws.Range("I4").Value2 = 0.1098676276771345
$ws.Range("J4").Value2 = 0.1098676276771345
$ws.Range("O4").Value2 = 0.2802068273109533
$ws.Range("P4").Value2 = 0.2802068273109533
$ws.Range("Q4").Value2 = 333.4918494144756
$ws.Range("R4").Value2 = 3001.42664473028
$ws.Range("S4").Value2 = 0.03078565937559094
$ws.Range("T4").Value2 = 0.03078565937559094
$ws.Range("H5").Value2 = 88.285005
$ws.Range("I5").Value2 = 0.1493784449296822
$ws.Range("J5").Value2 = 0.1493784449296822
$ws.Range("M5").Value2 = 3.759736666666667
$ws.Range("N5").Value2 = 11.27921
$ws.Range("O5").Value2 = 0.0683751702595819
$ws.Range("P5").Value2 = 0.06837517025958188
$ws.Range("Q5").Value2 = 110.64279013845
$ws.Range("R5").Value2 = 995.78511124605
$ws.Range("S5").Value2 = 0.0102137766051786
$ws.Range("T5").Value2 = 0.0102137766051786
$ws.Range("H6").Value2 = 88.285005
$ws.Range("I6").Value2 = 0.1493784449296822
$ws.Range("J6").Value2 = 0.1493784449296822
$ws.Range("O6").Value2 = 0.6514180024294648
$ws.Range("P6").Value2 = 0.6514180024294647
$ws.Range("R6").Value2 = 9486.957700496565
$ws.Range("S6").Value2 = 0.09730780820211336
$ws.Range("T6").Value2 = 0.09730780820211336
$ws.Range("H7").Value2 = 88.285005
$ws.Range("I7").Value2 = 0.1493784449296822
$ws.Range("J7").Value2 = 0.1493784449296822
$ws.Range("O7").Value2 = 0.2802068273109533
$ws.Range("P7").Value2 = 0.2802068273109533
$ws.Range("R7").Value2 = 4080.80573176545
$ws.Range("S7").Value2 = 0.0418568601223902
$ws.Range("T7").Value2 = 0.0418568601223902
$ws.Range("I8").Value2 = 0.7407539273931834
$ws.Range("J8").Value2 = 0.7407539273931834
$ws.Range("M8").Value2 = 3.759736666666667
$ws.Range("N8").Value2 = 11.27921
$ws.Range("O8").Value2 = 0.0683751702595819
$ws.Range("P8").Value2 = 0.06837517025958188
$ws.Range("Q8").Value2 = 548.6673888684389
$ws.Range("R8").Value2 = 4938.00649981595
$ws.Range("S8").Value2 = 0.05064917590596288
$ws.Range("T8").Value2 = 0.05064917590596286
$ws.Range("I9").Value2 = 0.7407539273931834
$ws.Range("J9").Value2 = 0.7407539273931834
$ws.Range("O9").Value2 = 0.6514180024294648
$ws.Range("P9").Value2 = 0.6514180024294647
$ws.Range("S9").Value2 = 0.4825404436742483
$ws.Range("T9").Value2 = 0.4825404436742482
$ws.Range("I10").Value2 = 0.7407539273931834
$ws.Range("J10").Value2 = 0.7407539273931834
$ws.Range("O10").Value2 = 0.2802068273109533
$ws.Range("P10").Value2 = 0.2802068273109533
$ws.Range("S10").Value2 = 0.2075643078129722
$ws.Range("T10").Value2 = 0.2075643078129721
